# "added test for default value"
# Adds a third worksheet (Sheet3) to the workbook, populated with job data
# plus two new date columns (END_DATE / START_DATE) and a couple of
# deliberately-blank cells (row 2 / row 19) that exercise "default value"
# handling, then makes Sheet3 the active sheet/tab.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet3: new worksheet, inserted right after Sheet2 -------------------
$ws3 = $wb.Worksheets.Add([Type]::Missing, $ws2)
$ws3.Name = "Sheet3"

    $ws3.Range("A1").Value = "JOB_ID"
    $ws3.Range("B1").Value = "JOB_TITLE"
    $ws3.Range("C1").Value = "MIN_SALARY"
    $ws3.Range("D1").Value = "MAX_SALARY"
    $ws3.Range("E1").Value = "END_DATE"
    $ws3.Range("F1").Value = "START_DATE"
    $ws3.Range("D2").Value = 40000
    $ws3.Range("A3").Value = "AD_VP"
    $ws3.Range("B3").Value = "Administration Vice President"
    $ws3.Range("C3").Value = 2008
    $ws3.Range("D3").Value = 30000
    $ws3.Range("E3").Value = "15.03.05"
    $ws3.Range("F3").Value = "21.09.97"
    $ws3.Range("A4").Value = "AD_ASST"
    $ws3.Range("B4").Value = "Administration Assistant"
    $ws3.Range("C4").Value = 2500
    $ws3.Range("D4").Value = 6000
    $ws3.Range("E4").Value = "15.03.05"
    $ws3.Range("F4").Value = "28.10.01"
    $ws3.Range("A5").Value = "FI_MGR"
    $ws3.Range("B5").Value = "Finance Manager"
    $ws3.Range("C5").Value = 8200
    $ws3.Range("D5").Value = 16000
    $ws3.Range("E5").Value = "19.12.07"
    $ws3.Range("F5").Value = "17.02.04"
    $ws3.Range("A6").Value = "FI_ACCOUNT"
    $ws3.Range("B6").Value = "Accountant"
    $ws3.Range("C6").Value = 4200
    $ws3.Range("D6").Value = 9000
    $ws3.Range("E6").Value = "31.12.07"
    $ws3.Range("F6").Value = "24.03.06"
    $ws3.Range("A7").Value = "AC_MGR"
    $ws3.Range("B7").Value = "Accounting Manager"
    $ws3.Range("C7").Value = 8200
    $ws3.Range("D7").Value = 16000
    $ws3.Range("E7").Value = "31.12.07"
    $ws3.Range("F7").Value = "'01.01.07"
    $ws3.Range("A8").Value = "AC_ACCOUNT"
    $ws3.Range("B8").Value = "Public Accountant"
    $ws3.Range("C8").Value = 4200
    $ws3.Range("D8").Value = 9000
    $ws3.Range("E8").Value = "17.06.01"
    $ws3.Range("F8").Value = "17.09.95"
    $ws3.Range("A9").Value = "SA_MAN"
    $ws3.Range("B9").Value = "Sales Manager"
    $ws3.Range("C9").Value = 10000
    $ws3.Range("D9").Value = 20080
    $ws3.Range("E9").Value = "31.12.06"
    $ws3.Range("F9").Value = "24.03.06"
    $ws3.Range("A10").Value = "SA_REP"
    $ws3.Range("B10").Value = "Sales Representative"
    $ws3.Range("C10").Value = 6000
    $ws3.Range("D10").Value = 12008
    $ws3.Range("E10").Value = "31.12.07"
    $ws3.Range("F10").Value = "'01.01.07"
    $ws3.Range("A11").Value = "PU_MAN"
    $ws3.Range("B11").Value = "Purchasing Manager"
    $ws3.Range("C11").Value = 8000
    $ws3.Range("D11").Value = 15000
    $ws3.Range("E11").Value = "31.12.06"
    $ws3.Range("F11").Value = "'01.07.02"
    $ws3.Range("A12").Value = "PU_CLERK"
    $ws3.Range("B12").Value = "Purchasing Clerk"
    $ws3.Range("C12").Value = 2500
    $ws3.Range("D12").Value = 5500
    $ws3.Range("E12").Value = "24.07.06"
    $ws3.Range("F12").Value = "21.09.97"
    $ws3.Range("A13").Value = "ST_MAN"
    $ws3.Range("B13").Value = "Stock Manager"
    $ws3.Range("C13").Value = 5500
    $ws3.Range("D13").Value = 8500
    $ws3.Range("E13").Value = "27.10.01"
    $ws3.Range("F13").Value = "28.10.01"
    $ws3.Range("A14").Value = "ST_CLERK"
    $ws3.Range("B14").Value = "Stock Clerk"
    $ws3.Range("C14").Value = 2008
    $ws3.Range("D14").Value = 5000
    $ws3.Range("E14").Value = "15.03.05"
    $ws3.Range("F14").Value = "17.02.04"
    $ws3.Range("A15").Value = "SH_CLERK"
    $ws3.Range("B15").Value = "Shipping Clerk"
    $ws3.Range("C15").Value = 2500
    $ws3.Range("D15").Value = 5500
    $ws3.Range("E15").Value = "19.12.07"
    $ws3.Range("F15").Value = "24.03.06"
    $ws3.Range("A16").Value = "IT_PROG"
    $ws3.Range("B16").Value = "Programmer"
    $ws3.Range("C16").Value = 4000
    $ws3.Range("D16").Value = 10000
    $ws3.Range("E16").Value = "31.12.07"
    $ws3.Range("F16").Value = "'01.01.07"
    $ws3.Range("A17").Value = "MK_MAN"
    $ws3.Range("B17").Value = "Marketing Manager"
    $ws3.Range("C17").Value = 9000
    $ws3.Range("D17").Value = 15000
    $ws3.Range("E17").Value = "31.12.07"
    $ws3.Range("F17").Value = "17.09.95"
    $ws3.Range("A18").Value = "MK_REP"
    $ws3.Range("B18").Value = "Marketing Representative"
    $ws3.Range("C18").Value = 4000
    $ws3.Range("D18").Value = 9000
    $ws3.Range("E18").Value = "17.06.01"
    $ws3.Range("F18").Value = "24.03.06"
    $ws3.Range("A19").Value = "HR_REP"
    $ws3.Range("A20").Value = "PR_REP"
    $ws3.Range("B20").Value = "Public Relations Representative"
    $ws3.Range("C20").Value = 4500
    $ws3.Range("D20").Value = 10500
    $ws3.Range("E20").Value = "31.12.07"
    $ws3.Range("F20").Value = "'01.07.02"
    $ws3.Range("D22").Value = "s"

# --- View/selection bookkeeping -------------------------------------------
# Sheet1: no longer the selected tab; zoom reset to 100%, selection -> D20
$ws1.Range("D20").Select()
$excel.ActiveWindow.Zoom = 100

# Sheet2: zoom reset to 100%, selection collapsed to the single cell D25
$ws2.Activate()
$ws2.Range("D25").Select()
$excel.ActiveWindow.Zoom = 100

# Sheet3: becomes the active/selected tab, selection on A19
$ws3.Activate()
$ws3.Range("A19").Select()
$excel.ActiveWindow.Zoom = 100
